$d = $word.ActiveDocument

function Get-ParagraphByText($needle) {
    $rng = $d.Content
    $found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find paragraph containing: $needle"
    }
    return $rng.Paragraphs(1)
}

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- 1) Add a trailing run containing a single space after the "ads" sentence ---
$adsPara = Get-ParagraphByText("DreamGrid does not use any ads, spyware, or other bad practices.")
$adsXml = $pkgHeader + `
    '<w:p w14:paraId="19CDF64B" w14:textId="77777777" w:rsidR="00A5108A" w:rsidRPr="00A5108A" w:rsidRDefault="00A5108A" w:rsidP="00853E93">' + `
        '<w:r w:rsidRPr="00A5108A"><w:t>DreamGrid does not use any ads, spyware, or other bad practices.</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '</w:p>' + $pkgFooter
$adsPara.Range.InsertXML($adsXml)

# --- 2) Drop "2011 " from the bold copyright line, keeping the trailing "." run intact ---
$copyrightPara = Get-ParagraphByText("DreamGrid is Copyright 2011 by Outworldz, LLC")
$copyrightXml = $pkgHeader + `
    '<w:p w14:paraId="199C383D" w14:textId="0F0D88F9" w:rsidR="00A5108A" w:rsidRDefault="00A5108A" w:rsidP="00A5108A">' + `
        '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' + `
        '<w:r w:rsidRPr="00A5108A"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>DreamGrid is Copyright  by Outworldz, LLC</w:t></w:r>' + `
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>.</w:t></w:r>' + `
    '</w:p>' + $pkgFooter
$copyrightPara.Range.InsertXML($copyrightXml)

# --- 3) Merge the split "Outworldz," / " LLC does not collect..." runs and drop the proofErr markers ---
$dataPara = Get-ParagraphByText("LLC does not collect any personally identifiable information")
$dataXml = $pkgHeader + `
    '<w:p w14:paraId="3BFBE402" w14:textId="77777777" w:rsidR="00A5108A" w:rsidRPr="00A5108A" w:rsidRDefault="00A5108A" w:rsidP="002741D2">' + `
        '<w:r w:rsidRPr="00A5108A"><w:t>Outworldz, LLC does not collect any personally identifiable information. No personal details of your site, such as Opensimulator passwords, usernames are other grid contents are sent to Outworldz, LLC. No Mesh, textures, inventory, or other items are sent to Outworldz, LLC.</w:t></w:r>' + `
    '</w:p>' + $pkgFooter
$dataPara.Range.InsertXML($dataXml)
